$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two MEXICO rows (rows 14-15); remaining rows shift up,
# so the former SOUTH KOREA rows (16-17) become rows 14-15.
$ws.Rows(14).Delete()
$ws.Rows(14).Delete()

# Odds/value corrections across the sheet
$ws.Range("C3").Value = "02:03"
$ws.Range("M3").Value = 1.07
$ws.Range("N3").Value = 9
$ws.Range("Q3").Value = 2.15
$ws.Range("R3").Value = 1.67
$ws.Range("G4").Value = 1.95
$ws.Range("G5").Value = 1.75
$ws.Range("I5").Value = 4.5
$ws.Range("J5").Value = 2.38
$ws.Range("M5").Value = 1.05
$ws.Range("N5").Value = 11
$ws.Range("U5").Value = 1.73
$ws.Range("V5").Value = 2
$ws.Range("Z5").Value = 15
$ws.Range("AB5").Value = 23
$ws.Range("AG5").Value = 201
$ws.Range("AW5").Value = 6
$ws.Range("G6").Value = 4.5
$ws.Range("I6").Value = 1.75
$ws.Range("J6").Value = 4.75
$ws.Range("L6").Value = 2.38
$ws.Range("Z6").Value = 51
$ws.Range("AF6").Value = 51
$ws.Range("AI6").Value = 8.5
$ws.Range("AR6").Value = 101
$ws.Range("G7").Value = 2.05
$ws.Range("G8").Value = 2.25
$ws.Range("H8").Value = 3.6
$ws.Range("AN8").Value = 4.5
$ws.Range("G11").Value = 3.1
$ws.Range("I11").Value = 2.35
$ws.Range("L11").Value = 3.2
$ws.Range("AE11").Value = 15
$ws.Range("AI11").Value = 11
$ws.Range("AJ11").Value = 10
$ws.Range("AK11").Value = 23
$ws.Range("AO11").Value = 17
$ws.Range("AQ11").Value = 51
$ws.Range("AR11").Value = 81
$ws.Range("AZ11").Value = 51
$ws.Range("Q13").Value = 2
$ws.Range("R13").Value = 1.85
$ws.Range("M14").Value = 1.05
$ws.Range("N14").Value = 11
$ws.Range("H15").Value = 3.2
$ws.Range("J15").Value = 3.2
$ws.Range("M15").Value = 1.07
$ws.Range("N15").Value = 9
$ws.Range("AO15").Value = 15
$ws.Range("AQ15").Value = 51
$ws.Range("AS15").Value = 201
